$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last refreshed" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 08:29"

# Row 6 - India: updated case numbers
$ws.Range("B6").Value = 1437976
$ws.Range("C6").Value = 1957
$ws.Range("D6").Value = 918906
$ws.Range("E6").Value = 486244
$ws.Range("G6").Value = 14
$ws.Range("H6").Value = 32826

# Row 54 - Afganistan: updated case numbers
$ws.Range("B54").Value = 36263
$ws.Range("C54").Value = 106
$ws.Range("D54").Value = 25198
$ws.Range("E54").Value = 9796
$ws.Range("G54").Value = 10
$ws.Range("H54").Value = 1269

# Row 56 now holds Kirguistan's updated data (Kirguistan overtook Ghana)
$ws.Range("A56").Value = "Kirguistan"
$ws.Range("B56").Value = 33296
$ws.Range("C56").Value = 483
$ws.Range("D56").Value = 21205
$ws.Range("E56").Value = 10790
$ws.Range("G56").Value = 24
$ws.Range("H56").Value = 1301

# Row 57 now holds Ghana's (previous) data
$ws.Range("A57").Value = "Ghana"
$ws.Range("B57").Value = 32969
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 29494
$ws.Range("E57").Value = 3307
$ws.Range("H57").Value = 168

# Row 75 - El Salvador: updated case numbers
$ws.Range("D75").Value = 7667
$ws.Range("E75").Value = 6555
$ws.Range("G75").Value = 8
$ws.Range("H75").Value = 408

# Row 109 - Tailandia: updated case numbers
$ws.Range("B109").Value = 3295
$ws.Range("C109").Value = 4
$ws.Range("D109").Value = 3111
$ws.Range("E109").Value = 126

# Row 143 now holds Georgia's updated data (Georgia overtook Niger)
$ws.Range("A143").Value = "Georgia"
$ws.Range("B143").Value = 1137
$ws.Range("C143").Value = 6
$ws.Range("D143").Value = 922
$ws.Range("E143").Value = 199
$ws.Range("H143").Value = 16

# Row 144 now holds Niger's (previous) data
$ws.Range("A144").Value = "Niger"
$ws.Range("B144").Value = 1136
$ws.Range("C144").Value = 0
$ws.Range("D144").Value = 1027
$ws.Range("E144").Value = 40
$ws.Range("H144").Value = 69

# Rows 210/211: Islas Malvinas and Groenlandia swap names (identical data)
$ws.Range("A210").Value = "Groenlandia"
$ws.Range("A211").Value = "Islas Malvinas"
